# Generate Report for Handoff
# Adds 3 new "Ready for handoff" rows (for a new source file
# 4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md plus its two dependent
# .png assets) to the Overview / zh-cn / de-de sheets, resizes their
# tables, and wires up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview" (sheet1) - summary rows 5..7, columns A-G
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ovRows = @(
    @{ A="4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md"; B="e2e\4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md"; C=".md";  E="Ready for handoff"; F="Ready for handoff"; G="2016-10-17 17:52:16";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/e2e/4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md" },
    @{ A="492e9433-7836-4fec-9297-feaa1fe6b8d0.png"; B="e2e\492e9433-7836-4fec-9297-feaa1fe6b8d0.png"; C=".png"; E="Ready for handoff"; F="Ready for handoff"; G="2016-10-17 17:52:16";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/e2e/492e9433-7836-4fec-9297-feaa1fe6b8d0.png" },
    @{ A="7a0756ee-120b-438d-bd5a-bf2c09c1f21b.png"; B="e2e\7a0756ee-120b-438d-bd5a-bf2c09c1f21b.png"; C=".png"; E="Ready for handoff"; F="Ready for handoff"; G="2016-10-17 17:52:16";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/e2e/7a0756ee-120b-438d-bd5a-bf2c09c1f21b.png" }
)

# (Source Path = column D is blank for every Overview row, so it is
# intentionally left untouched below - writing "" would clear the cell.)
$r = 5
foreach ($row in $ovRows) {
    $ov.Cells.Item($r, 1).Value = $row.A
    $ov.Cells.Item($r, 2).Value = $row.B
    $ov.Cells.Item($r, 3).Value = $row.C
    $ov.Cells.Item($r, 5).Value = $row.E
    $ov.Cells.Item($r, 6).Value = $row.F
    $ov.Cells.Item($r, 7).Value = $row.G
    $ov.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ov.Hyperlinks.Add($ov.Cells.Item($r, 2), $row.Url, "", "", $row.B) | Out-Null

    $r = $r + 1
}

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G7"))

# ---------------------------------------------------------------
# Sheet "zh-cn" (sheet2) / "de-de" (sheet3) - detail rows 5..7,
# columns A-P. Both sheets share the same shape/content except for
# column G (target xlf file name), column H (zh-cn uses the real
# handoff timestamp, de-de mirrors the "Ready for handoff" text).
# ---------------------------------------------------------------
function Fill-LocaleSheet($ws, $xlfSuffix, $colH) {
    $rows = @(
        @{ A="4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md"; B=".md";  G="4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.695c73b5fbd5f98899a1f8163b07de8a542e4e90.$xlfSuffix.xlf";
           K="0001-01-01 00:00:00"; M="True"; N="";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md" },
        @{ A="492e9433-7836-4fec-9297-feaa1fe6b8d0.png"; B=".png"; G="31bd0427d75129feee567f4dc0ea5ef5d1559416.png";
           K="0001-01-01 00:00:00"; M="True(Dependency)"; N="e2e\4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/492e9433-7836-4fec-9297-feaa1fe6b8d0.png" },
        @{ A="7a0756ee-120b-438d-bd5a-bf2c09c1f21b.png"; B=".png"; G="e10e53892d6bcb42d1092e448dbce17e620630b8.png";
           K="0001-01-01 00:00:00"; M="True(Dependency)"; N="e2e\4f8a8d95-ba7e-47dc-9f89-7b815cf86f80.md";
           Url="https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d83113c12ecdef18de1f31e10fcd148ce7002d4f/7a0756ee-120b-438d-bd5a-bf2c09c1f21b.png" }
    )

    # Columns I, J, L, P are blank ("") for every new row here - left
    # untouched on purpose since writing "" clears/omits the cell, which
    # renders identically to an explicit empty string.
    $r = 5
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row.A                    # A Source File Name
        $ws.Cells.Item($r, 2).Value = $row.B                    # B File Extension
        $ws.Cells.Item($r, 3).Value = "Ready for handoff"       # C Status
        $ws.Cells.Item($r, 4).Value = "e2e"                     # D Source Path
        $ws.Cells.Item($r, 5).Value = "ht"                      # E Priority
        $ws.Cells.Item($r, 6).Value = "False"                   # F Content Duplicate
        $ws.Cells.Item($r, 7).Value = $row.G                    # G Latest Handoff File
        $ws.Cells.Item($r, 8).Value = $colH                     # H Latest Handoff Datetime
        $ws.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Cells.Item($r, 11).Value = $row.K                   # K Latest Handback DateTime
        $ws.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Cells.Item($r, 13).Value = $row.M                   # M To be localized
        if ($row.N -ne "") {
            $ws.Cells.Item($r, 14).Value = $row.N               # N Dependency From
        }
        $ws.Cells.Item($r, 15).Value = "False"                  # O Has metadata

        $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), $row.Url, "", "", $row.A) | Out-Null

        $r = $r + 1
    }
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Fill-LocaleSheet $zhcn "zh-cn" "2016-10-17 17:51:54"
$zhcnTable = $zhcn.ListObjects.Item(1)
$zhcnTable.Resize($zhcn.Range("A1:P7"))

$dede = $wb.Worksheets.Item("de-de")
Fill-LocaleSheet $dede "de-de" "Ready for handoff"
$dedeTable = $dede.ListObjects.Item(1)
$dedeTable.Resize($dede.Range("A1:P7"))
